$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.157.74'
$ws.Range('D3').Value = '1.676.75'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'214.24"
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = "'22.88"
$ws.Range('E8').Value = '  +7.14%  '
$ws.Range('D9').Value = "'0.262"
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('D10').Value = "'0.0621"
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = "'0.0890"
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '1.913.56'
$ws.Range('D13').Value = '1.674.96'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('E14').Value = '  +2.30%  '
$ws.Range('D15').Value = "'0.562"
$ws.Range('E15').Value = '  +4.85%  '
$ws.Range('D16').Value = "'66.55"
$ws.Range('D17').Value = '27.128.34'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = "'235.53"
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '0.0₃0742'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('D20').Value = "'7.83"
$ws.Range('E20').Value = '  -4.16%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = "'4.55"
$ws.Range('E22').Value = '  +1.83%  '
$ws.Range('E23').Value = '  +3.03%  '
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('D25').Value = "'148.40"
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').Value = "'7.48"
$ws.Range('E26').Value = '  +2.50%  '
$ws.Range('D27').Value = "'16.44"
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('D28').Value = "'0.113"
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = "'0.0500"
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '1.547.69'
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('E34').Value = '  +1.85%  '
$ws.Range('E35').Value = '  -3.89%  '
$ws.Range('D36').Value = "'0.609"
$ws.Range('E36').Value = '  +3.57%  '
$ws.Range('E37').Value = '  +3.12%  '
$ws.Range('D38').Value = "'2.38"
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('D41').Value = "'70.04"
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('E42').Value = '  +4.49%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').Value = '1.822.36'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = "'0.781"
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').Value = "'1.65"
$ws.Range('E47').Value = '  +6.27%  '
$ws.Range('D48').Value = "'89.77"
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('D50').Value = "'8.25"
$ws.Range('E50').Value = '  +2.73%  '
$ws.Range('E51').Value = '  +0.18%  '
